# Weekly price-sheet refresh for "Hortaliza, Mapocho Venta Directa de
# Santiago - Sandia": dates and the associated volume/price/origin/
# quality figures for several entries are updated to reflect the new
# reporting week (commit: "Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44497
$ws.Range("J2").Value = 250
$ws.Range("D5").Value = 44504
$ws.Range("J5").Value = 200
$ws.Range("K5").Value = 800
$ws.Range("L5").Value = 800
$ws.Range("M5").Value = 800
$ws.Range("N5").Value = "$/kilo (volumen en unidades)"
$ws.Range("P5").Value = 800
$ws.Range("D6").Value = 44223
$ws.Range("H6").Value = "Americana O Klondike"
$ws.Range("I6").Value = "Extra"
$ws.Range("J6").Value = 340
$ws.Range("K6").Value = 2500
$ws.Range("L6").Value = 2500
$ws.Range("M6").Value = 2500
$ws.Range("N6").Value = "$/unidad"
$ws.Range("O6").Value = "Región de O'Higgins"
$ws.Range("P6").Value = 2500
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 400
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = 2000
$ws.Range("P7").Value = 2000
$ws.Range("I8").Value = "Segunda"
$ws.Range("J8").Value = 300
$ws.Range("K8").Value = 1500
$ws.Range("L8").Value = 1500
$ws.Range("M8").Value = 1500
$ws.Range("P8").Value = 1500
$ws.Range("I9").Value = "Tercera"
$ws.Range("J9").Value = 160
$ws.Range("K9").Value = 1000
$ws.Range("L9").Value = 1000
$ws.Range("M9").Value = 1000
$ws.Range("P9").Value = 1000
$ws.Range("D10").Value = 44495
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 200
$ws.Range("K10").Value = 800
$ws.Range("L10").Value = 800
$ws.Range("M10").Value = 800
$ws.Range("N10").Value = "$/kilo (volumen en unidades)"
$ws.Range("O10").Value = "Perú"
$ws.Range("P10").Value = 800
$ws.Range("D11").Value = 44510
$ws.Range("J11").Value = 250
$ws.Range("D12").Value = 44477
$ws.Range("J12").Value = 80
$ws.Range("D13").Value = 44167
$ws.Range("J13").Value = 400
$ws.Range("K13").Value = 5000
$ws.Range("L13").Value = 5000
$ws.Range("M13").Value = 5000
$ws.Range("N13").Value = "$/unidad"
$ws.Range("O13").Value = "Región de O'Higgins"
$ws.Range("P13").Value = 5000
$ws.Range("D14").Value = 44167
$ws.Range("I14").Value = "Segunda"
$ws.Range("J14").Value = 560
$ws.Range("K14").Value = 3000
$ws.Range("L14").Value = 3000
$ws.Range("M14").Value = 3000
$ws.Range("O14").Value = "Región de O'Higgins"
$ws.Range("P14").Value = 3000
$ws.Range("D15").Value = 44167
$ws.Range("I15").Value = "Tercera"
$ws.Range("J15").Value = 450
$ws.Range("K15").Value = 2000
$ws.Range("L15").Value = 2000
$ws.Range("M15").Value = 2000
$ws.Range("N15").Value = "$/unidad"
$ws.Range("O15").Value = "Región de O'Higgins"
$ws.Range("P15").Value = 2000
$ws.Range("D16").Value = 44483
$ws.Range("J16").Value = 120
$ws.Range("K16").Value = 800
$ws.Range("L16").Value = 800
$ws.Range("M16").Value = 800
$ws.Range("N16").Value = "$/kilo (volumen en unidades)"
$ws.Range("O16").Value = "Perú"
$ws.Range("P16").Value = 800
$ws.Range("D17").Value = 44488
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 150
$ws.Range("K17").Value = 800
$ws.Range("L17").Value = 800
$ws.Range("M17").Value = 800
$ws.Range("N17").Value = "$/kilo (volumen en unidades)"
$ws.Range("O17").Value = "Perú"
$ws.Range("P17").Value = 800
$ws.Range("D18").Value = 44305
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 2500
$ws.Range("L18").Value = 2500
$ws.Range("M18").Value = 2500
$ws.Range("O18").Value = "Perú"
$ws.Range("P18").Value = 2500
$ws.Range("D19").Value = 44194
$ws.Range("I19").Value = "Extra"
$ws.Range("J19").Value = 120
$ws.Range("K19").Value = 3500
$ws.Range("L19").Value = 3500
$ws.Range("M19").Value = 3500
$ws.Range("N19").Value = "$/unidad"
$ws.Range("O19").Value = "Región de O'Higgins"
$ws.Range("P19").Value = 3500
$ws.Range("D20").Value = 44194
$ws.Range("J20").Value = 200
$ws.Range("K20").Value = 3000
$ws.Range("L20").Value = 3000
$ws.Range("M20").Value = 3000
$ws.Range("N20").Value = "$/unidad"
$ws.Range("O20").Value = "Región de O'Higgins"
$ws.Range("P20").Value = 3000
$ws.Range("D21").Value = 44491
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 150
$ws.Range("K21").Value = 800
$ws.Range("L21").Value = 800
$ws.Range("M21").Value = 800
$ws.Range("N21").Value = "$/kilo (volumen en unidades)"
$ws.Range("O21").Value = "Perú"
$ws.Range("P21").Value = 800
$ws.Range("D22").Value = 44312
$ws.Range("J22").Value = 180
$ws.Range("K22").Value = 2500
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = 2500
$ws.Range("O22").Value = "Perú"
$ws.Range("P22").Value = 2500

